$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 holds a date-like label ("11-08-2021") that must be stored as TEXT
# (shared string), not auto-converted to a date serial. Force text format,
# assign, then clear the format so the cell ends up with no explicit style
# (matching row 2's plain cells).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "11-08-2021"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = 350000
$ws.Range("C3").Value = 435000
$ws.Range("D3").Value = 350000
$ws.Range("E3").Value = 335000
$ws.Range("F3").Value = 15000
$ws.Range("G3").Value = 2.8
